$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header + timestamp values for the new 'time_taken' column (F),
# matching the source data's per-row extraction timestamps.
$values = New-Object 'object[,]' 113,1
$values[0,0] = "time_taken"
$values[1,0] = "2021-10-05 13:39:11.705093"
$values[2,0] = "2021-10-05 13:39:11.705104"
$values[3,0] = "2021-10-05 13:39:11.705107"
$values[4,0] = "2021-10-05 13:39:11.705110"
$values[5,0] = "2021-10-05 13:39:11.705113"
$values[6,0] = "2021-10-05 13:39:11.705115"
$values[7,0] = "2021-10-05 13:39:11.705118"
$values[8,0] = "2021-10-05 13:39:11.705120"
$values[9,0] = "2021-10-05 13:39:11.705123"
$values[10,0] = "2021-10-05 13:39:11.705125"
$values[11,0] = "2021-10-05 13:39:11.705128"
$values[12,0] = "2021-10-05 13:39:11.705130"
$values[13,0] = "2021-10-05 13:39:11.705132"
$values[14,0] = "2021-10-05 13:39:11.705135"
$values[15,0] = "2021-10-05 13:39:11.705137"
$values[16,0] = "2021-10-05 13:39:11.705140"
$values[17,0] = "2021-10-05 13:39:11.705142"
$values[18,0] = "2021-10-05 13:39:11.705145"
$values[19,0] = "2021-10-05 13:39:11.705148"
$values[20,0] = "2021-10-05 13:39:11.705150"
$values[21,0] = "2021-10-05 13:39:11.705152"
$values[22,0] = "2021-10-05 13:39:11.705155"
$values[23,0] = "2021-10-05 13:39:11.705157"
$values[24,0] = "2021-10-05 13:39:11.705160"
$values[25,0] = "2021-10-05 13:39:11.705162"
$values[26,0] = "2021-10-05 13:39:11.705165"
$values[27,0] = "2021-10-05 13:39:11.705167"
$values[28,0] = "2021-10-05 13:39:11.705170"
$values[29,0] = "2021-10-05 13:39:11.705172"
$values[30,0] = "2021-10-05 13:39:11.705175"
$values[31,0] = "2021-10-05 13:39:11.705177"
$values[32,0] = "2021-10-05 13:39:11.705179"
$values[33,0] = "2021-10-05 13:39:11.705182"
$values[34,0] = "2021-10-05 13:39:11.705185"
$values[35,0] = "2021-10-05 13:39:11.705187"
$values[36,0] = "2021-10-05 13:39:11.705189"
$values[37,0] = "2021-10-05 13:39:11.705192"
$values[38,0] = "2021-10-05 13:39:11.705194"
$values[39,0] = "2021-10-05 13:39:11.705197"
$values[40,0] = "2021-10-05 13:39:11.705199"
$values[41,0] = "2021-10-05 13:39:11.705202"
$values[42,0] = "2021-10-05 13:39:11.705205"
$values[43,0] = "2021-10-05 13:39:11.705207"
$values[44,0] = "2021-10-05 13:39:11.705210"
$values[45,0] = "2021-10-05 13:39:11.705212"
$values[46,0] = "2021-10-05 13:39:11.705215"
$values[47,0] = "2021-10-05 13:39:11.705217"
$values[48,0] = "2021-10-05 13:39:11.705220"
$values[49,0] = "2021-10-05 13:39:11.705222"
$values[50,0] = "2021-10-05 13:39:11.705224"
$values[51,0] = "2021-10-05 13:39:11.705227"
$values[52,0] = "2021-10-05 13:39:11.705229"
$values[53,0] = "2021-10-05 13:39:11.705232"
$values[54,0] = "2021-10-05 13:39:11.705235"
$values[55,0] = "2021-10-05 13:39:11.705237"
$values[56,0] = "2021-10-05 13:39:11.705239"
$values[57,0] = "2021-10-05 13:39:11.705242"
$values[58,0] = "2021-10-05 13:39:11.705244"
$values[59,0] = "2021-10-05 13:39:11.705247"
$values[60,0] = "2021-10-05 13:39:11.705249"
$values[61,0] = "2021-10-05 13:39:11.705251"
$values[62,0] = "2021-10-05 13:39:11.705254"
$values[63,0] = "2021-10-05 13:39:11.705256"
$values[64,0] = "2021-10-05 13:39:11.705259"
$values[65,0] = "2021-10-05 13:39:11.705262"
$values[66,0] = "2021-10-05 13:39:11.705265"
$values[67,0] = "2021-10-05 13:39:11.705267"
$values[68,0] = "2021-10-05 13:39:11.705270"
$values[69,0] = "2021-10-05 13:39:11.705272"
$values[70,0] = "2021-10-05 13:39:11.705275"
$values[71,0] = "2021-10-05 13:39:11.705277"
$values[72,0] = "2021-10-05 13:39:11.705279"
$values[73,0] = "2021-10-05 13:39:11.705282"
$values[74,0] = "2021-10-05 13:39:11.705284"
$values[75,0] = "2021-10-05 13:39:11.705287"
$values[76,0] = "2021-10-05 13:39:11.705289"
$values[77,0] = "2021-10-05 13:39:11.705293"
$values[78,0] = "2021-10-05 13:39:11.705297"
$values[79,0] = "2021-10-05 13:39:11.705299"
$values[80,0] = "2021-10-05 13:39:11.705302"
$values[81,0] = "2021-10-05 13:39:11.705304"
$values[82,0] = "2021-10-05 13:39:11.705307"
$values[83,0] = "2021-10-05 13:39:11.705309"
$values[84,0] = "2021-10-05 13:39:11.705312"
$values[85,0] = "2021-10-05 13:39:11.705314"
$values[86,0] = "2021-10-05 13:39:11.705317"
$values[87,0] = "2021-10-05 13:39:11.705319"
$values[88,0] = "2021-10-05 13:39:11.705321"
$values[89,0] = "2021-10-05 13:39:11.705324"
$values[90,0] = "2021-10-05 13:39:11.705326"
$values[91,0] = "2021-10-05 13:39:11.705329"
$values[92,0] = "2021-10-05 13:39:11.705331"
$values[93,0] = "2021-10-05 13:39:11.705335"
$values[94,0] = "2021-10-05 13:39:11.705338"
$values[95,0] = "2021-10-05 13:39:11.705340"
$values[96,0] = "2021-10-05 13:39:11.705343"
$values[97,0] = "2021-10-05 13:39:11.705345"
$values[98,0] = "2021-10-05 13:39:11.705348"
$values[99,0] = "2021-10-05 13:39:11.705350"
$values[100,0] = "2021-10-05 13:39:11.705353"
$values[101,0] = "2021-10-05 13:39:11.705355"
$values[102,0] = "2021-10-05 13:39:11.705358"
$values[103,0] = "2021-10-05 13:39:11.705360"
$values[104,0] = "2021-10-05 13:39:11.705363"
$values[105,0] = "2021-10-05 13:39:11.705365"
$values[106,0] = "2021-10-05 13:39:11.705368"
$values[107,0] = "2021-10-05 13:39:11.705370"
$values[108,0] = "2021-10-05 13:39:11.705373"
$values[109,0] = "2021-10-05 13:39:11.705377"
$values[110,0] = "2021-10-05 13:39:11.705380"
$values[111,0] = "2021-10-05 13:39:11.705383"
$values[112,0] = "2021-10-05 13:39:11.705385"

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -ne 113) { $lastRow = 113 }
$targetRange = $ws.Range("F1:F" + $lastRow)
$targetRange.Value = $values

# Match header cell F1's style to the other header cells (bold, bordered, centered)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "Applied time_taken column to F1:F$lastRow"
